# Updates computed profit/price figures (columns H-N) for specific leve rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, reflecting
# refreshed market-board data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3040
$ws.Range("I64").Value = 2642.6667
$ws.Range("J64").Value = 3323.8096
$ws.Range("K64").Value = 2642.6667
$ws.Range("L64").Value = 3323.8096
$ws.Range("M64").Value = -2394.6667
$ws.Range("N64").Value = -3819.8096
$ws.Range("H67").Value = 3040
$ws.Range("I67").Value = 2642.6667
$ws.Range("J67").Value = 3323.8096
$ws.Range("K67").Value = 2642.6667
$ws.Range("L67").Value = 3323.8096
$ws.Range("M67").Value = -1784.6667
$ws.Range("N67").Value = -5039.809600000001
$ws.Range("H70").Value = 63742.062
$ws.Range("I70").Value = 251070
$ws.Range("J70").Value = 1299.4166
$ws.Range("K70").Value = 753210
$ws.Range("L70").Value = 3898.2498
$ws.Range("M70").Value = -752940
$ws.Range("N70").Value = -4438.2498
$ws.Range("H73").Value = 63742.062
$ws.Range("I73").Value = 251070
$ws.Range("J73").Value = 1299.4166
$ws.Range("K73").Value = 753210
$ws.Range("L73").Value = 3898.2498
$ws.Range("M73").Value = -752274
$ws.Range("N73").Value = -5770.2498
$ws.Range("H116").Value = 12502437
$ws.Range("I116").Value = 16669017
$ws.Range("J116").Value = 2697.25
$ws.Range("K116").Value = 16669017
$ws.Range("L116").Value = 2697.25
$ws.Range("M116").Value = -16665575
$ws.Range("N116").Value = -9581.25
$ws.Range("H131").Value = 4149.552
$ws.Range("I131").Value = 770.5
$ws.Range("J131").Value = 5670.125
$ws.Range("K131").Value = 2311.5
$ws.Range("L131").Value = 17010.375
$ws.Range("M131").Value = 2728.5
$ws.Range("N131").Value = -27090.375
$ws.Range("H135").Value = 948.80646
$ws.Range("I135").Value = 892.2083
$ws.Range("J135").Value = 1142.8572
$ws.Range("K135").Value = 8029.8747
$ws.Range("L135").Value = 10285.7148
$ws.Range("M135").Value = -5494.8747
$ws.Range("N135").Value = -15355.7148
$ws.Range("H136").Value = 66195
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 66195
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 66195
$ws.Range("N136").Value = -76395
$ws.Range("H137").Value = 1236.289
$ws.Range("I137").Value = 1018.125
$ws.Range("J137").Value = 2981.6
$ws.Range("K137").Value = 3054.375
$ws.Range("L137").Value = 8944.799999999999
$ws.Range("M137").Value = -504.375
$ws.Range("N137").Value = -14044.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 22633.334
$ws.Range("I22").Value = 9000
$ws.Range("J22").Value = 36266.668
$ws.Range("K22").Value = 9000
$ws.Range("L22").Value = 36266.668
$ws.Range("M22").Value = -8701
$ws.Range("N22").Value = -36864.668
$ws.Range("H45").Value = 1302.1538
$ws.Range("I45").Value = 1177.091
$ws.Range("J45").Value = 1990
$ws.Range("K45").Value = 1177.091
$ws.Range("L45").Value = 1990
$ws.Range("M45").Value = -800.0909999999999
$ws.Range("N45").Value = -2744
$ws.Range("H132").Value = 2160.7
$ws.Range("I132").Value = 1334.1482
$ws.Range("J132").Value = 3877.3845
$ws.Range("K132").Value = 4002.4446
$ws.Range("L132").Value = 11632.1535
$ws.Range("M132").Value = -1472.4446
$ws.Range("N132").Value = -16692.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 5546
$ws.Range("I22").Value = 6151.1113
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 6151.1113
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = -5978.1113
$ws.Range("N22").Value = -446

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9994.727999999999
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 9994.727999999999
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 9994.727999999999
$ws.Range("N4").Value = -10218.728
$ws.Range("H58").Value = 994.03705
$ws.Range("I58").Value = 1044.5
$ws.Range("J58").Value = 772
$ws.Range("K58").Value = 1044.5
$ws.Range("L58").Value = 772
$ws.Range("M58").Value = -841.5
$ws.Range("N58").Value = -1178
$ws.Range("H122").Value = 2220.7144
$ws.Range("I122").Value = 2004.2693
$ws.Range("J122").Value = 2846
$ws.Range("K122").Value = 6012.8079
$ws.Range("L122").Value = 8538
$ws.Range("M122").Value = -3562.8079
$ws.Range("N122").Value = -13438
$ws.Range("H134").Value = 2583.1538
$ws.Range("I134").Value = 2433.3333
$ws.Range("J134").Value = 2711.5715
$ws.Range("K134").Value = 7299.999899999999
$ws.Range("L134").Value = 8134.7145
$ws.Range("M134").Value = -4764.999899999999
$ws.Range("N134").Value = -13204.7145
$ws.Range("H136").Value = 994.03705
$ws.Range("I136").Value = 1044.5
$ws.Range("J136").Value = 772
$ws.Range("K136").Value = 3133.5
$ws.Range("L136").Value = 2316
$ws.Range("M136").Value = -583.5
$ws.Range("N136").Value = -7416

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 262.94116
$ws.Range("I11").Value = 83.375
$ws.Range("J11").Value = 422.55554
$ws.Range("K11").Value = 250.125
$ws.Range("L11").Value = 1267.66662
$ws.Range("M11").Value = -110.125
$ws.Range("N11").Value = -1547.66662
$ws.Range("H21").Value = 4169477.2
$ws.Range("I21").Value = 703.3333
$ws.Range("J21").Value = 5956094.5
$ws.Range("K21").Value = 2109.9999
$ws.Range("L21").Value = 17868283.5
$ws.Range("M21").Value = -1936.9999
$ws.Range("N21").Value = -17868629.5
$ws.Range("H22").Value = 100000350
$ws.Range("I22").Value = 100000350
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 300001050
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -300000881
$ws.Range("N22").Value = ""
$ws.Range("H25").Value = 1357.6666
$ws.Range("I25").Value = 144
$ws.Range("J25").Value = 2874.75
$ws.Range("K25").Value = 432
$ws.Range("L25").Value = 8624.25
$ws.Range("M25").Value = -263
$ws.Range("N25").Value = -8962.25
$ws.Range("H27").Value = 100000350
$ws.Range("I27").Value = 100000350
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 300001050
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -300000948
$ws.Range("N27").Value = ""
$ws.Range("H30").Value = 1357.6666
$ws.Range("I30").Value = 144
$ws.Range("J30").Value = 2874.75
$ws.Range("K30").Value = 432
$ws.Range("L30").Value = 8624.25
$ws.Range("M30").Value = -330
$ws.Range("N30").Value = -8828.25
$ws.Range("H34").Value = 552.625
$ws.Range("I34").Value = 197.33333
$ws.Range("J34").Value = 634.61536
$ws.Range("K34").Value = 591.99999
$ws.Range("L34").Value = 1903.84608
$ws.Range("M34").Value = -507.99999
$ws.Range("N34").Value = -2071.84608
$ws.Range("H39").Value = 2519.75
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2519.75
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 7559.25
$ws.Range("N39").Value = -8147.25
$ws.Range("H40").Value = 404.2857
$ws.Range("I40").Value = 157.5
$ws.Range("J40").Value = 733.3333
$ws.Range("K40").Value = 630
$ws.Range("L40").Value = 2933.3332
$ws.Range("M40").Value = -561
$ws.Range("N40").Value = -3071.3332
$ws.Range("H55").Value = 4419.1665
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 4589.5654
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 13768.6962
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -14122.6962
$ws.Range("H64").Value = 5123.923
$ws.Range("I64").Value = 2970.6667
$ws.Range("J64").Value = 5769.9
$ws.Range("K64").Value = 8912.000100000001
$ws.Range("L64").Value = 17309.7
$ws.Range("M64").Value = -8642.000100000001
$ws.Range("N64").Value = -17849.7
$ws.Range("H67").Value = 5123.923
$ws.Range("I67").Value = 2970.6667
$ws.Range("J67").Value = 5769.9
$ws.Range("K67").Value = 8912.000100000001
$ws.Range("L67").Value = 17309.7
$ws.Range("M67").Value = -7976.000100000001
$ws.Range("N67").Value = -19181.7
$ws.Range("H68").Value = 1773.8422
$ws.Range("I68").Value = 1450
$ws.Range("J68").Value = 1811.9412
$ws.Range("K68").Value = 4350
$ws.Range("L68").Value = 5435.8236
$ws.Range("M68").Value = -3539
$ws.Range("N68").Value = -7057.8236
$ws.Range("H71").Value = 1773.8422
$ws.Range("I71").Value = 1450
$ws.Range("J71").Value = 1811.9412
$ws.Range("K71").Value = 13050
$ws.Range("L71").Value = 16307.4708
$ws.Range("M71").Value = -8994
$ws.Range("N71").Value = -24419.4708
$ws.Range("H94").Value = 3472.5
$ws.Range("I94").Value = 1980
$ws.Range("J94").Value = 3970
$ws.Range("K94").Value = 5940
$ws.Range("L94").Value = 11910
$ws.Range("M94").Value = -5264
$ws.Range("N94").Value = -13262
$ws.Range("H98").Value = 2219
$ws.Range("I98").Value = 2757.25
$ws.Range("J98").Value = 1501.3334
$ws.Range("K98").Value = 8271.75
$ws.Range("L98").Value = 4504.0002
$ws.Range("M98").Value = -6773.75
$ws.Range("N98").Value = -7500.0002
$ws.Range("H112").Value = 3510.2856
$ws.Range("I112").Value = 1513.5
$ws.Range("J112").Value = 4309
$ws.Range("K112").Value = 4540.5
$ws.Range("L112").Value = 12927
$ws.Range("M112").Value = -3432.5
$ws.Range("N112").Value = -15143
$ws.Range("H133").Value = 5274
$ws.Range("I133").Value = 2343.3333
$ws.Range("J133").Value = 6530
$ws.Range("K133").Value = 7029.999899999999
$ws.Range("L133").Value = 19590
$ws.Range("M133").Value = -1969.999899999999
$ws.Range("N133").Value = -29710

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 1250
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1250
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1250
$ws.Range("N17").Value = -1586
$ws.Range("H132").Value = 2086.261
$ws.Range("I132").Value = 1755.1111
$ws.Range("J132").Value = 3278.4
$ws.Range("K132").Value = 5265.3333
$ws.Range("L132").Value = 9835.200000000001
$ws.Range("M132").Value = -2735.3333
$ws.Range("N132").Value = -14895.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9928
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 9928
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 9928
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = -10152

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4003
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4003
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4003
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -5251
$ws.Range("H65").Value = 4003
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4003
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 20015
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -26255
$ws.Range("H122").Value = 18385178
$ws.Range("I122").Value = 19233308
$ws.Range("J122").Value = 15628754
$ws.Range("K122").Value = 57699924
$ws.Range("L122").Value = 46886262
$ws.Range("M122").Value = -57697474
$ws.Range("N122").Value = -46891162
$ws.Range("H126").Value = 16897
$ws.Range("I126").Value = 16897
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 50691
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -48221
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = ""
